$wb = $excel.ActiveWorkbook

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H5").Value = 2999.5
$ws.Range("I5").Value = 2999.5
$ws.Range("K5").Value = 2999.5
$ws.Range("M5").Value = -2884.5
$ws.Range("H32").Value = 19600.8
$ws.Range("I32").Value = 0
$ws.Range("J32").Value = 19600.8
$ws.Range("K32").Value = 0
$ws.Range("L32").Value = 19600.8
$ws.Range("N32").Value = -20252.8
$ws.Range("H49").Value = 1635.5
$ws.Range("I49").Value = 1349.25
$ws.Range("J49").Value = 1921.75
$ws.Range("K49").Value = 4047.75
$ws.Range("L49").Value = 5765.25
$ws.Range("M49").Value = -3911.75
$ws.Range("N49").Value = -6037.25
$ws.Range("H113").Value = 2466.375
$ws.Range("I113").Value = 1654.1666
$ws.Range("J113").Value = 4903
$ws.Range("K113").Value = 1654.1666
$ws.Range("L113").Value = 4903
$ws.Range("M113").Value = 1599.8334
$ws.Range("N113").Value = -11411
$ws.Range("H127").Value = 3000
$ws.Range("I127").Value = 2950
$ws.Range("K127").Value = 8850
$ws.Range("M127").Value = -3890
$ws.Range("H133").Value = 115000
$ws.Range("J133").Value = 115000
$ws.Range("L133").Value = 115000
$ws.Range("N133").Value = -125120
$ws.Range("H137").Value = 3034.2222
$ws.Range("I137").Value = 2786.4285
$ws.Range("K137").Value = 8359.2855
$ws.Range("M137").Value = -5809.2855
$ws.Range("H141").Value = 4700
$ws.Range("J141").Value = 4500
$ws.Range("L141").Value = 13500
$ws.Range("N141").Value = -23860
$ws.Range("M32").ClearContents()

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 76933900
$ws.Range("I61").Value = 90920020
$ws.Range("J61").Value = 10250.5
$ws.Range("K61").Value = 90920020
$ws.Range("L61").Value = 10250.5
$ws.Range("M61").Value = -90919808
$ws.Range("N61").Value = -10674.5
$ws.Range("H74").Value = 125011224
$ws.Range("I74").Value = 125011224
$ws.Range("J74").Value = 0
$ws.Range("K74").Value = 125011224
$ws.Range("L74").Value = 0
$ws.Range("M74").Value = -125010350
$ws.Range("H77").Value = 125011224
$ws.Range("I77").Value = 125011224
$ws.Range("J77").Value = 0
$ws.Range("K77").Value = 625056120
$ws.Range("L77").Value = 0
$ws.Range("M77").Value = -625051752
$ws.Range("H103").Value = 99999
$ws.Range("J103").Value = 99999
$ws.Range("L103").Value = 99999
$ws.Range("N103").Value = -102343
$ws.Range("H104").Value = 26333.334
$ws.Range("I104").Value = 39000
$ws.Range("J104").Value = 20000
$ws.Range("K104").Value = 39000
$ws.Range("L104").Value = 20000
$ws.Range("M104").Value = -35506
$ws.Range("N104").Value = -26988
$ws.Range("H132").Value = 6671073
$ws.Range("I132").Value = 6671073
$ws.Range("K132").Value = 20013219
$ws.Range("M132").Value = -20010689
$ws.Range("H136").Value = 76933900
$ws.Range("I136").Value = 90920020
$ws.Range("J136").Value = 10250.5
$ws.Range("K136").Value = 272760060
$ws.Range("L136").Value = 30751.5
$ws.Range("M136").Value = -272757510
$ws.Range("N136").Value = -35851.5
$ws.Range("H138").Value = 0
$ws.Range("J138").Value = 0
$ws.Range("L138").Value = 0
$ws.Range("N74").ClearContents()
$ws.Range("N77").ClearContents()
$ws.Range("N138").ClearContents()

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H134").Value = 25002768
$ws.Range("I134").Value = 26317652
$ws.Range("K134").Value = 78952956
$ws.Range("M134").Value = -78950421

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 8584.333000000001
$ws.Range("I7").Value = 25073.25
$ws.Range("J7").Value = 339.875
$ws.Range("K7").Value = 25073.25
$ws.Range("L7").Value = 339.875
$ws.Range("M7").Value = -24960.25
$ws.Range("N7").Value = -565.875
$ws.Range("H31").Value = 9817.119000000001
$ws.Range("I31").Value = 7373.091
$ws.Range("J31").Value = 12505.55
$ws.Range("K31").Value = 7373.091
$ws.Range("L31").Value = 12505.55
$ws.Range("M31").Value = -7078.091
$ws.Range("N31").Value = -13095.55
$ws.Range("H34").Value = 9817.119000000001
$ws.Range("I34").Value = 7373.091
$ws.Range("J34").Value = 12505.55
$ws.Range("K34").Value = 7373.091
$ws.Range("L34").Value = 12505.55
$ws.Range("M34").Value = -7171.091
$ws.Range("N34").Value = -12909.55
$ws.Range("H50").Value = 29990
$ws.Range("J50").Value = 29990
$ws.Range("L50").Value = 29990
$ws.Range("N50").Value = -31240
$ws.Range("H58").Value = 29420152
$ws.Range("I58").Value = 38470970
$ws.Range("K58").Value = 38470970
$ws.Range("M58").Value = -38470767
$ws.Range("H99").Value = 1959.6
$ws.Range("I99").Value = 1700
$ws.Range("J99").Value = 2998
$ws.Range("K99").Value = 1700
$ws.Range("L99").Value = 2998
$ws.Range("M99").Value = -202
$ws.Range("N99").Value = -5994
$ws.Range("H126").Value = 1959.6
$ws.Range("I126").Value = 1700
$ws.Range("J126").Value = 2998
$ws.Range("K126").Value = 5100
$ws.Range("L126").Value = 8994
$ws.Range("M126").Value = -2630
$ws.Range("N126").Value = -13934
$ws.Range("H132").Value = 25642382
$ws.Range("I132").Value = 25642382
$ws.Range("K132").Value = 76927146
$ws.Range("M132").Value = -76924616
$ws.Range("H136").Value = 29420152
$ws.Range("I136").Value = 38470970
$ws.Range("K136").Value = 115412910
$ws.Range("M136").Value = -115410360

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 1940624.6
$ws.Range("J4").Value = 2225610.8
$ws.Range("L4").Value = 6676832.399999999
$ws.Range("N4").Value = -6677056.399999999
$ws.Range("H92").Value = 1299.5
$ws.Range("I92").Value = 600
$ws.Range("J92").Value = 1999
$ws.Range("K92").Value = 1800
$ws.Range("L92").Value = 5997
$ws.Range("M92").Value = -552
$ws.Range("N92").Value = -8493
$ws.Range("H98").Value = 928.3333
$ws.Range("J98").Value = 598.75
$ws.Range("L98").Value = 1796.25
$ws.Range("N98").Value = -4792.25
$ws.Range("H141").Value = 555
$ws.Range("I141").Value = 555
$ws.Range("K141").Value = 1665
$ws.Range("M141").Value = 3515

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H3").Value = 1000
$ws.Range("I3").Value = 1000
$ws.Range("K3").Value = 1000
$ws.Range("M3").Value = -884
$ws.Range("H7").Value = 1447569.9
$ws.Range("I7").Value = 2017598
$ws.Range("J7").Value = 22499.5
$ws.Range("K7").Value = 2017598
$ws.Range("L7").Value = 22499.5
$ws.Range("M7").Value = -2017486
$ws.Range("N7").Value = -22723.5
$ws.Range("H8").Value = 1447569.9
$ws.Range("I8").Value = 2017598
$ws.Range("J8").Value = 22499.5
$ws.Range("K8").Value = 2017598
$ws.Range("L8").Value = 22499.5
$ws.Range("M8").Value = -2017459
$ws.Range("N8").Value = -22777.5
$ws.Range("H11").Value = 8026.8887
$ws.Range("I11").Value = 9576.933999999999
$ws.Range("K11").Value = 9576.933999999999
$ws.Range("M11").Value = -9437.933999999999
$ws.Range("H18").Value = 0
$ws.Range("I18").Value = 0
$ws.Range("K18").Value = 0
$ws.Range("H24").Value = 39332
$ws.Range("J24").Value = 39998
$ws.Range("L24").Value = 39998
$ws.Range("N24").Value = -40344
$ws.Range("H132").Value = 7817207
$ws.Range("I132").Value = 8336280
$ws.Range("K132").Value = 25008840
$ws.Range("M132").Value = -25006310
$ws.Range("M18").ClearContents()

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H11").Value = 7000
$ws.Range("J11").Value = 7000
$ws.Range("L11").Value = 7000
$ws.Range("N11").Value = -7280
$ws.Range("H13").Value = 13666.333
$ws.Range("J13").Value = 14499.5
$ws.Range("L13").Value = 14499.5
$ws.Range("N13").Value = -14779.5
$ws.Range("H17").Value = 20000
$ws.Range("I17").Value = 20000
$ws.Range("J17").Value = 0
$ws.Range("K17").Value = 20000
$ws.Range("L17").Value = 0
$ws.Range("M17").Value = -19830
$ws.Range("H20").Value = 18000
$ws.Range("I20").Value = 0
$ws.Range("J20").Value = 18000
$ws.Range("K20").Value = 0
$ws.Range("L20").Value = 18000
$ws.Range("N20").Value = -18452
$ws.Range("H55").Value = 203.23334
$ws.Range("I55").Value = 144.66667
$ws.Range("J55").Value = 261.8
$ws.Range("K55").Value = 144.66667
$ws.Range("L55").Value = 261.8
$ws.Range("M55").Value = 28.33332999999999
$ws.Range("N55").Value = -607.8
$ws.Range("H122").Value = 4770.3706
$ws.Range("I122").Value = 4899.7827
$ws.Range("K122").Value = 14699.3481
$ws.Range("M122").Value = -12249.3481
$ws.Range("H132").Value = 17152176
$ws.Range("I132").Value = 17787372
$ws.Range("J132").Value = 1900
$ws.Range("K132").Value = 53362116
$ws.Range("L132").Value = 5700
$ws.Range("M132").Value = -53359586
$ws.Range("N132").Value = -10760
$ws.Range("N17").ClearContents()
$ws.Range("M20").ClearContents()

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H15").Value = 12000
$ws.Range("I15").Value = 12000
$ws.Range("K15").Value = 12000
$ws.Range("M15").Value = -11712
$ws.Range("H17").Value = 6451
$ws.Range("I17").Value = 4676.5
$ws.Range("J17").Value = 10000
$ws.Range("K17").Value = 4676.5
$ws.Range("L17").Value = 10000
$ws.Range("M17").Value = -4504.5
$ws.Range("N17").Value = -10344
$ws.Range("H113").Value = 441.16666
$ws.Range("I113").Value = 368.8
$ws.Range("J113").Value = 492.85715
$ws.Range("K113").Value = 1106.4
$ws.Range("L113").Value = 1478.57145
$ws.Range("M113").Value = 1063.6
$ws.Range("N113").Value = -5818.571449999999
$ws.Range("H132").Value = 19237376
$ws.Range("I132").Value = 33336060
$ws.Range("K132").Value = 100008180
$ws.Range("M132").Value = -100005650
$ws.Range("H136").Value = 35714972
$ws.Range("I136").Value = 38461850
$ws.Range("K136").Value = 115385550
$ws.Range("M136").Value = -115383000
$ws.Range("H141").Value = 66664.664
$ws.Range("I141").Value = 0
$ws.Range("J141").Value = 0
$ws.Range("K141").Value = 0
$ws.Range("L141").Value = 66664.664
$ws.Range("N141").Value = -77024.664
$ws.Range("M141").ClearContents()
